$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force just the Price (D) cells we are about to rewrite to stay text,
# so numeric-looking values like "612.58" or "7.80" are not coerced
# into floats/doubles by Excel (they were inlineStr/text in the source).
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D6").NumberFormat = "@"
$ws.Range("D13:D17").NumberFormat = "@"
$ws.Range("D19:D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37:D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46:D47").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.075.45"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.698.69"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "612.58"
$ws.Range("D6").Value = "158.72"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  +4.66%  "
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "0.0000209"
$ws.Range("E13").Value = "  +9.59%  "
$ws.Range("D14").Value = "30.22"
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("D15").Value = "3.182.42"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "65.934.73"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "2.695.42"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "4.91"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "7.80"
$ws.Range("D21").Value = "359.65"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "71.40"
$ws.Range("E22").Value = "  +3.21%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  +17.26%  "
$ws.Range("D25").Value = "9.98"
$ws.Range("E25").Value = "  +5.70%  "
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("D29").Value = "8.31"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "535.98"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "6.71"
$ws.Range("E34").Value = "  +4.54%  "
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "20.80"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").Value = "163.08"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "168.65"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").Value = "23.86"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "2.32"
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("D49").Value = "0.659"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").Value = "20.96"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("D51").Value = "0.0997"
$ws.Range("E51").Value = "  +1.42%  "
